# Auto-generated Excel COM-interop script applying the Asura_Profits market-data refresh.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across the ALC, ARM,
# BSM, CRP, CUL, GSM, LTW and WVR sheets with freshly pulled market values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 72070.14
$ws.Range("I107").Value = 100618.3
$ws.Range("J107").Value = 699.75
$ws.Range("K107").Value = 100618.3
$ws.Range("L107").Value = 699.75
$ws.Range("M107").Value = -98698.3
$ws.Range("N107").Value = -4539.75

$ws.Range("H129").Value = 1271.125
$ws.Range("J129").Value = 1309.5172
$ws.Range("L129").Value = 3928.5516
$ws.Range("N129").Value = -13928.5516

$ws.Range("H132").Value = 1534.2941
$ws.Range("I132").Value = 1376.3507
$ws.Range("J132").Value = 3054.5
$ws.Range("K132").Value = 4129.0521
$ws.Range("L132").Value = 9163.5
$ws.Range("M132").Value = -1599.0521
$ws.Range("N132").Value = -14223.5

$ws.Range("H133").Value = 72896.25
$ws.Range("J133").Value = 72896.25
$ws.Range("L133").Value = 72896.25
$ws.Range("N133").Value = -83016.25

$ws.Range("H134").Value = 113465.414
$ws.Range("J134").Value = 113465.414
$ws.Range("L134").Value = 113465.414
$ws.Range("N134").Value = -123605.414

$ws.Range("H137").Value = 1908.6
$ws.Range("I137").Value = 1804
$ws.Range("J137").Value = 2136.818
$ws.Range("K137").Value = 5412
$ws.Range("L137").Value = 6410.454000000001
$ws.Range("M137").Value = -2862
$ws.Range("N137").Value = -11510.454

$ws.Range("H138").Value = 2252679
$ws.Range("I138").Value = 5559570.5
$ws.Range("J138").Value = 6488.4907
$ws.Range("K138").Value = 16678711.5
$ws.Range("L138").Value = 19465.4721
$ws.Range("M138").Value = -16673571.5
$ws.Range("N138").Value = -29745.4721

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 251500
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 251500
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -251726

$ws.Range("H32").Value = 38820.24
$ws.Range("I32").Value = 32650.162
$ws.Range("J32").Value = 66144.86
$ws.Range("K32").Value = 32650.162
$ws.Range("L32").Value = 66144.86
$ws.Range("M32").Value = -32363.162
$ws.Range("N32").Value = -66718.86

$ws.Range("H61").Value = 2869.7097
$ws.Range("I61").Value = 2355.682
$ws.Range("J61").Value = 4126.222
$ws.Range("K61").Value = 2355.682
$ws.Range("L61").Value = 4126.222
$ws.Range("M61").Value = -2143.682
$ws.Range("N61").Value = -4550.222

$ws.Range("H74").Value = 1491.9412
$ws.Range("I74").Value = 1518.68
$ws.Range("J74").Value = 1417.6666
$ws.Range("K74").Value = 1518.68
$ws.Range("L74").Value = 1417.6666
$ws.Range("M74").Value = -644.6800000000001
$ws.Range("N74").Value = -3165.6666

$ws.Range("H77").Value = 1491.9412
$ws.Range("I77").Value = 1518.68
$ws.Range("J77").Value = 1417.6666
$ws.Range("K77").Value = 7593.400000000001
$ws.Range("L77").Value = 7088.333000000001
$ws.Range("M77").Value = -3225.400000000001
$ws.Range("N77").Value = -15824.333

$ws.Range("H102").Value = 84858.336
$ws.Range("I102").Value = 1500
$ws.Range("K102").Value = 1500
$ws.Range("M102").Value = 122

$ws.Range("H116").Value = 251500
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 251500
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -256088

$ws.Range("H132").Value = 2098.1633
$ws.Range("I132").Value = 1767.9231
$ws.Range("J132").Value = 3386.1
$ws.Range("K132").Value = 5303.7693
$ws.Range("L132").Value = 10158.3
$ws.Range("M132").Value = -2773.7693
$ws.Range("N132").Value = -15218.3

$ws.Range("H136").Value = 2869.7097
$ws.Range("I136").Value = 2355.682
$ws.Range("J136").Value = 4126.222
$ws.Range("K136").Value = 7067.045999999999
$ws.Range("L136").Value = 12378.666
$ws.Range("M136").Value = -4517.045999999999
$ws.Range("N136").Value = -17478.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 251500
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 251500
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -251728

$ws.Range("H64").Value = 261.64285
$ws.Range("I64").Value = 334
$ws.Range("J64").Value = 131.4
$ws.Range("K64").Value = 334
$ws.Range("L64").Value = 131.4
$ws.Range("M64").Value = -109
$ws.Range("N64").Value = -581.4

$ws.Range("H67").Value = 261.64285
$ws.Range("I67").Value = 334
$ws.Range("J67").Value = 131.4
$ws.Range("K67").Value = 334
$ws.Range("L67").Value = 131.4
$ws.Range("M67").Value = 446
$ws.Range("N67").Value = -1691.4

$ws.Range("H105").Value = 2222.0645
$ws.Range("I105").Value = 2229.4666
$ws.Range("K105").Value = 2229.4666
$ws.Range("M105").Value = -482.4666000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1702.931
$ws.Range("I134").Value = 1526.0834
$ws.Range("J134").Value = 2551.8
$ws.Range("K134").Value = 4578.2502
$ws.Range("L134").Value = 7655.400000000001
$ws.Range("M134").Value = -2043.2502
$ws.Range("N134").Value = -12725.4

$ws.Range("H138").Value = 80734
$ws.Range("J138").Value = 80734
$ws.Range("L138").Value = 80734
$ws.Range("N138").Value = -91014

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6850
$ws.Range("J80").Value = 3542.8572
$ws.Range("L80").Value = 10628.5716
$ws.Range("N80").Value = -12500.5716

$ws.Range("H83").Value = 6850
$ws.Range("J83").Value = 3542.8572
$ws.Range("L83").Value = 31885.7148
$ws.Range("N83").Value = -41245.7148

$ws.Range("H122").Value = 7953.375
$ws.Range("I122").Value = 595.5
$ws.Range("K122").Value = 5359.5
$ws.Range("M122").Value = -2909.5

$ws.Range("H124").Value = 2502.5
$ws.Range("I124").Value = 755
$ws.Range("K124").Value = 2265
$ws.Range("M124").Value = 2645

$ws.Range("H125").Value = 3018.111
$ws.Range("I125").Value = 2022.5
$ws.Range("J125").Value = 3814.6
$ws.Range("K125").Value = 6067.5
$ws.Range("L125").Value = 11443.8
$ws.Range("M125").Value = -1147.5
$ws.Range("N125").Value = -21283.8

$ws.Range("H129").Value = 2002077.4
$ws.Range("I129").Value = 950
$ws.Range("J129").Value = 2176088.5
$ws.Range("K129").Value = 2850
$ws.Range("L129").Value = 6528265.5
$ws.Range("M129").Value = 2150
$ws.Range("N129").Value = -6538265.5

$ws.Range("H131").Value = 47626096
$ws.Range("J131").Value = 58826220
$ws.Range("L131").Value = 176478660
$ws.Range("N131").Value = -176488740

$ws.Range("H136").Value = 4888.75
$ws.Range("I136").Value = 1420
$ws.Range("J136").Value = 6045
$ws.Range("K136").Value = 4260
$ws.Range("L136").Value = 18135
$ws.Range("M136").Value = 840
$ws.Range("N136").Value = -28335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 60000000
$ws.Range("I3").Value = 10000000
$ws.Range("J3").Value = 110000000
$ws.Range("K3").Value = 10000000
$ws.Range("L3").Value = 110000000
$ws.Range("M3").Value = -9999884
$ws.Range("N3").Value = -110000232

$ws.Range("H102").Value = 2039.0526
$ws.Range("I102").Value = 1994.625
$ws.Range("K102").Value = 1994.625
$ws.Range("M102").Value = -372.625

$ws.Range("H107").Value = 6843.067
$ws.Range("I107").Value = 8457.166999999999
$ws.Range("J107").Value = 386.66666
$ws.Range("K107").Value = 8457.166999999999
$ws.Range("L107").Value = 386.66666
$ws.Range("M107").Value = -6537.166999999999
$ws.Range("N107").Value = -4226.66666

$ws.Range("H122").Value = 5042.857
$ws.Range("I122").Value = 4216.6665
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 12649.9995
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -10199.9995
$ws.Range("N122").Value = -34900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 41666.668
$ws.Range("I11").Value = 50000
$ws.Range("J11").Value = 37500
$ws.Range("K11").Value = 50000
$ws.Range("L11").Value = 37500
$ws.Range("M11").Value = -49860
$ws.Range("N11").Value = -37780

$ws.Range("H46").Value = 1306.5834
$ws.Range("I46").Value = 1005
$ws.Range("J46").Value = 1457.375
$ws.Range("K46").Value = 1005
$ws.Range("L46").Value = 1457.375
$ws.Range("M46").Value = -817
$ws.Range("N46").Value = -1833.375

$ws.Range("H103").Value = 90000
$ws.Range("J103").Value = 90000
$ws.Range("L103").Value = 90000
$ws.Range("N103").Value = -92344

$ws.Range("H134").Value = 38660
$ws.Range("J134").Value = 38660
$ws.Range("L134").Value = 38660
$ws.Range("N134").Value = -48800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 250011410
$ws.Range("J101").Value = 250011410
$ws.Range("L101").Value = 250011410
$ws.Range("N101").Value = -250017900
